$wb = $excel.ActiveWorkbook

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 433.26666
$ws.Range("I98").Value = 417
$ws.Range("K98").Value = 417
$ws.Range("M98").Value = 1081

# Sheet ALC, row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 500
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

# Sheet ALC, row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1053.8
$ws.Range("I101").Value = 299.66666
$ws.Range("J101").Value = 2185
$ws.Range("K101").Value = 898.9999799999999
$ws.Range("L101").Value = 6555
$ws.Range("M101").Value = 723.0000200000001
$ws.Range("N101").Value = -9799

# Sheet ALC, row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3500
$ws.Range("I106").Value = 4000
$ws.Range("K106").Value = 4000
$ws.Range("M106").Value = -3369

# Sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 492.6875
$ws.Range("I107").Value = 563.61536
$ws.Range("J107").Value = 185.33333
$ws.Range("K107").Value = 563.61536
$ws.Range("L107").Value = 185.33333
$ws.Range("M107").Value = 1356.38464
$ws.Range("N107").Value = -4025.33333

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4516.5
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 433.26666
$ws.Range("I122").Value = 417
$ws.Range("K122").Value = 1251
$ws.Range("M122").Value = 1199

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5958131
$ws.Range("I132").Value = 6255787.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 18767362.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -18764832.5
$ws.Range("N132").Value = -20060

# Sheet ALC, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2854.8333
$ws.Range("I135").Value = 717.5
$ws.Range("J135").Value = 4992.1665
$ws.Range("K135").Value = 6457.5
$ws.Range("L135").Value = 44929.4985
$ws.Range("M135").Value = -3922.5
$ws.Range("N135").Value = -49999.4985

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1955.8286
$ws.Range("I138").Value = 1731.0667
$ws.Range("J138").Value = 2124.4
$ws.Range("K138").Value = 5193.2001
$ws.Range("L138").Value = 6373.200000000001
$ws.Range("M138").Value = -53.20010000000002
$ws.Range("N138").Value = -16653.2

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2611.1177
$ws.Range("I61").Value = 1964.8334
$ws.Range("J61").Value = 2963.6365
$ws.Range("K61").Value = 1964.8334
$ws.Range("L61").Value = 2963.6365
$ws.Range("M61").Value = -1752.8334
$ws.Range("N61").Value = -3387.6365

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 864.28
$ws.Range("I74").Value = 819.5
$ws.Range("K74").Value = 819.5
$ws.Range("M74").Value = 54.5

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 864.28
$ws.Range("I77").Value = 819.5
$ws.Range("K77").Value = 4097.5
$ws.Range("M77").Value = 270.5

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1445.2
$ws.Range("I122").Value = 1279.0476
$ws.Range("J122").Value = 2317.5
$ws.Range("K122").Value = 3837.142800000001
$ws.Range("L122").Value = 6952.5
$ws.Range("M122").Value = -1387.142800000001
$ws.Range("N122").Value = -11852.5

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2611.1177
$ws.Range("I136").Value = 1964.8334
$ws.Range("J136").Value = 2963.6365
$ws.Range("K136").Value = 5894.5002
$ws.Range("L136").Value = 8890.9095
$ws.Range("M136").Value = -3344.5002
$ws.Range("N136").Value = -13990.9095

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2244.375
$ws.Range("I134").Value = 1923.3077
$ws.Range("J134").Value = 3635.6667
$ws.Range("K134").Value = 5769.9231
$ws.Range("L134").Value = 10907.0001
$ws.Range("M134").Value = -3234.9231
$ws.Range("N134").Value = -15977.0001

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2357.9832
$ws.Range("I31").Value = 1565.25
$ws.Range("J31").Value = 2901.5715
$ws.Range("K31").Value = 1565.25
$ws.Range("L31").Value = 2901.5715
$ws.Range("M31").Value = -1270.25
$ws.Range("N31").Value = -3491.5715

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2357.9832
$ws.Range("I34").Value = 1565.25
$ws.Range("J34").Value = 2901.5715
$ws.Range("K34").Value = 1565.25
$ws.Range("L34").Value = 2901.5715
$ws.Range("M34").Value = -1363.25
$ws.Range("N34").Value = -3305.5715

# Sheet CRP, row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 25822.818
$ws.Range("J74").Value = 25822.818
$ws.Range("L74").Value = 25822.818
$ws.Range("N74").Value = -27570.818

# Sheet CRP, row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 25822.818
$ws.Range("J77").Value = 25822.818
$ws.Range("L77").Value = 77468.454
$ws.Range("N77").Value = -86204.454

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3331
$ws.Range("I86").Value = 3044.4443
$ws.Range("K86").Value = 3044.4443
$ws.Range("M86").Value = -1921.4443

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3331
$ws.Range("I89").Value = 3044.4443
$ws.Range("K89").Value = 15222.2215
$ws.Range("M89").Value = -9606.2215

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14156.223
$ws.Range("I99").Value = 4032
$ws.Range("J99").Value = 34404.668
$ws.Range("K99").Value = 4032
$ws.Range("L99").Value = 34404.668
$ws.Range("M99").Value = -2534
$ws.Range("N99").Value = -37400.668

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1021.5714
$ws.Range("I107").Value = 1024.25
$ws.Range("J107").Value = 1018
$ws.Range("K107").Value = 1024.25
$ws.Range("L107").Value = 1018
$ws.Range("M107").Value = 895.75
$ws.Range("N107").Value = -4858

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 617.5
$ws.Range("I122").Value = 585
$ws.Range("J122").Value = 715
$ws.Range("K122").Value = 1755
$ws.Range("L122").Value = 2145
$ws.Range("M122").Value = 695
$ws.Range("N122").Value = -7045

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 14156.223
$ws.Range("I126").Value = 4032
$ws.Range("J126").Value = 34404.668
$ws.Range("K126").Value = 12096
$ws.Range("L126").Value = 103214.004
$ws.Range("M126").Value = -9626
$ws.Range("N126").Value = -108154.004

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5324.1904
$ws.Range("I132").Value = 7923.4443
$ws.Range("J132").Value = 3374.75
$ws.Range("K132").Value = 23770.3329
$ws.Range("L132").Value = 10124.25
$ws.Range("M132").Value = -21240.3329
$ws.Range("N132").Value = -15184.25

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1047.3334
$ws.Range("J68").Value = 1056.8
$ws.Range("L68").Value = 3170.4
$ws.Range("N68").Value = -4792.4

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1047.3334
$ws.Range("J71").Value = 1056.8
$ws.Range("L71").Value = 9511.199999999999
$ws.Range("N71").Value = -17623.2

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 618
$ws.Range("J122").Value = 650
$ws.Range("L122").Value = 5850
$ws.Range("N122").Value = -10750

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 765.12
$ws.Range("J131").Value = 775.6326
$ws.Range("L131").Value = 2326.8978
$ws.Range("N131").Value = -12406.8978

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1535.375
$ws.Range("I132").Value = 700.5
$ws.Range("J132").Value = 2370.25
$ws.Range("K132").Value = 6304.5
$ws.Range("L132").Value = 21332.25
$ws.Range("M132").Value = -3774.5
$ws.Range("N132").Value = -26392.25

# Sheet CUL, row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2165.2964
$ws.Range("I137").Value = 1672.5
$ws.Range("J137").Value = 2559.5334
$ws.Range("K137").Value = 5017.5
$ws.Range("L137").Value = 7678.600199999999
$ws.Range("M137").Value = 82.5
$ws.Range("N137").Value = -17878.6002

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 972.5
$ws.Range("I122").Value = 963.3333
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2889.9999
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -439.9998999999998
$ws.Range("N122").Value = -7900

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2016.5814
$ws.Range("I132").Value = 1836.4706
$ws.Range("J132").Value = 2697
$ws.Range("K132").Value = 5509.4118
$ws.Range("L132").Value = 8091
$ws.Range("M132").Value = -2979.4118
$ws.Range("N132").Value = -13151

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 73649.92999999999
$ws.Range("I40").Value = 251424.75
$ws.Range("J40").Value = 2540
$ws.Range("K40").Value = 251424.75
$ws.Range("L40").Value = 2540
$ws.Range("M40").Value = -251288.75
$ws.Range("N40").Value = -2812

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1888.7059
$ws.Range("I61").Value = 1867.8889
$ws.Range("J61").Value = 1912.125
$ws.Range("K61").Value = 1867.8889
$ws.Range("L61").Value = 1912.125
$ws.Range("M61").Value = -1665.8889
$ws.Range("N61").Value = -2316.125

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1888.7059
$ws.Range("I113").Value = 1867.8889
$ws.Range("J113").Value = 1912.125
$ws.Range("K113").Value = 1867.8889
$ws.Range("L113").Value = 1912.125
$ws.Range("M113").Value = 302.1111000000001
$ws.Range("N113").Value = -6252.125

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 561.25
$ws.Range("I113").Value = 390.61905
$ws.Range("J113").Value = 887
$ws.Range("K113").Value = 1171.85715
$ws.Range("L113").Value = 2661
$ws.Range("M113").Value = 998.14285
$ws.Range("N113").Value = -7001

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2435.0588
$ws.Range("I122").Value = 1674.75
$ws.Range("J122").Value = 3110.889
$ws.Range("K122").Value = 5024.25
$ws.Range("L122").Value = 9332.667000000001
$ws.Range("M122").Value = -2574.25
$ws.Range("N122").Value = -14232.667

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1312.079
$ws.Range("I136").Value = 481.0606
$ws.Range("K136").Value = 1443.1818
$ws.Range("M136").Value = 1106.8182
